# Updated elec ces modeling
# - RQSD-BRQSD: mark "nuclear" and "hydro" as qualifying RPS sources (0 -> 1)
# - RQSD-RQSD: break the live links to RQSD-BRQSD; this sheet now carries its
#   own independent (hard-coded) boolean values instead of formulas that
#   pulled from the 'RQSD-BRQSD' sheet. Only the "lignite" row keeps a
#   formula, now pointing at this sheet's own B2 (hard coal) instead of
#   'RQSD-BRQSD'!B13. "biomass" and "municipal solid waste" are flipped to 0.

$wb = $excel.ActiveWorkbook

$wsBRQSD = $wb.Worksheets.Item("RQSD-BRQSD")
$wsRQSD  = $wb.Worksheets.Item("RQSD-RQSD")

# --- RQSD-BRQSD (BAU) sheet: flip nuclear & hydro to qualify ---------------
$wsBRQSD.Range("B4").Value = 1
$wsBRQSD.Range("B5").Value = 1
$wsBRQSD.Range("B4:B5").Style = "Normal"

# --- RQSD-RQSD sheet: replace cross-sheet formulas with plain values -------
$wsRQSD.Range("B2").Value = 0
$wsRQSD.Range("B3").Value = 0
$wsRQSD.Range("B4").Value = 0
$wsRQSD.Range("B5").Value = 0
$wsRQSD.Range("B6").Value = 1
$wsRQSD.Range("B7").Value = 1
$wsRQSD.Range("B8").Value = 1
$wsRQSD.Range("B9").Value = 0
$wsRQSD.Range("B10").Value = 1
$wsRQSD.Range("B11").Value = 0
$wsRQSD.Range("B12").Value = 0
$wsRQSD.Range("B13").Formula = "=B2"
$wsRQSD.Range("B14").Value = 1
$wsRQSD.Range("B15").Value = 0
$wsRQSD.Range("B16").Value = 0
$wsRQSD.Range("B17").Value = 0

$wsRQSD.Range("B2:B17").Style = "Normal"

# --- Window/selection bookkeeping to match the saved view state ------------
$wsRQSD.Activate()
$wsRQSD.Range("B2").Select()

$wsBRQSD.Activate()
$wsBRQSD.Range("B18").Select()
